$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Add the new "Egg Timer Counter" test-data block FIRST, in the exact
#    order that makes the shared-string table grow: timerCount (A22) before
#    testTimerCounterValue (A21), so they land at uniqueCount indexes 12/13
#    respectively (matching the target sharedStrings.xml ordering), ahead of
#    the relocated hyperlink strings.
# ---------------------------------------------------------------------------

# Row 22 ("timerCount" field-name row, no special style) -> new shared string 12
$ws.Range("A22").Value() = "timerCount"

# Row 21 ("testTimerCounterValue" header row, shaded like the other headers)
# -> new shared string 13
$ws.Range("A17:C17").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)
$ws.Range("A21").Value() = "testTimerCounterValue"

# Row 23 (A23 numeric value 25, with new left/top/wrap-text style; B23 blank)
$ws.Range("A23").Value() = 25
$ws.Range("A23").HorizontalAlignment = -4131   # xlLeft
$ws.Range("A23").VerticalAlignment = -4160     # xlTop
$ws.Range("A23").WrapText = $true

# Row 24 ("endOfTestData" footer, reuses shared string 0)
$ws.Range("A24").Value() = "endOfTestData"

# ---------------------------------------------------------------------------
# 2) Replace the old "xyz@yahoo.comUsr" / "abc@123$" hyperlinked test data
#    (rows 3, 7, 11) with the new "myplayer@yahoo.comUsr" / "Rajuway@123$"
#    values, and strip their hyperlink styling (s="3") back to the default.
# ---------------------------------------------------------------------------

$ws.Range("A3:B3").ClearFormats()
$ws.Range("A3").Value() = "myplayer@yahoo.comUsr"
$ws.Range("B3").Value() = "Rajuway@123$"

$ws.Range("A7:B7").ClearFormats()
$ws.Range("A7").Value() = "myplayer@yahoo.comUsr"
$ws.Range("B7").Value() = "Rajuway@123$"

$ws.Range("A11:B11").ClearFormats()
$ws.Range("A11").Value() = "myplayer@yahoo.comUsr"
$ws.Range("B11").Value() = "Rajuway@123$"
$ws.Range("C11").Clear()

# ---------------------------------------------------------------------------
# 3) Drop the now-stale hyperlink definitions entirely.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 4) Update the saved selection to reflect where editing left off.
# ---------------------------------------------------------------------------
[void]$ws.Range("D12").Select()
